$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 13 de Mayo de 2020 a las 11:35"

# Row 18 (Belgica) updated stats
$ws.Range("B18").Value = 53981
$ws.Range("C18").Value = 202
$ws.Range("D18").Value = 13937
$ws.Range("E18").Value = 31201
$ws.Range("F18").Value = 420
$ws.Range("G18").Value = 82
$ws.Range("H18").Value = 8843

# Rows 173/174: swap Malaui / Polinesia Francesa order + update stats
# Row 173 becomes Malaui with updated data
$ws.Range("A173").Value = "Malaui"
$ws.Range("B173").Value = 63
$ws.Range("C173").Value = 6
$ws.Range("D173").Value = 24
$ws.Range("E173").Value = 36
$ws.Range("F173").Value = 1
$ws.Range("G173").Value = 0
$ws.Range("H173").Value = 3

# Row 174 becomes Polinesia Francesa with its (unchanged) data
$ws.Range("A174").Value = "Polinesia Francesa"
$ws.Range("B174").Value = 60
$ws.Range("C174").Value = 0
$ws.Range("D174").Value = 58
$ws.Range("E174").Value = 2
$ws.Range("F174").Value = 1
$ws.Range("G174").Value = 0
$ws.Range("H174").Value = 0
